$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Cells.Item(14, 6).Value = "Param2           <value>"
$ws.Cells.Item(16, 6).Value = "Power             <value>W"
$ws.Cells.Item(18, 6).Value = "Current           <value>A"
$ws.Cells.Item(20, 6).Value = "Voltage           <value>V"
$ws.Cells.Item(78, 6).Value = "Frequency       <value>Hz"
$ws.Cells.Item(80, 6).Value = "Param1           <value>"
$ws.Cells.Item(86, 2).Value = "SingleUseId195"
$ws.Cells.Item(86, 4).Value = "Center"
$ws.Cells.Item(86, 6).Value = "Enable Parameter"
$ws.Cells.Item(87, 2).Value = "SingleUseId196"
$ws.Cells.Item(87, 6).Value = "Main Menu"
$ws.Cells.Item(88, 2).Value = "SingleUseId197"
$ws.Cells.Item(88, 4).Value = "Left"
$ws.Cells.Item(88, 6).Value = "CPU Usage: <value>%"
$ws.Cells.Item(89, 2).Value = "SingleUseId198"
$ws.Cells.Item(89, 6).Value = "0"
$ws.Cells.Item(90, 2).Value = "SingleUseId199"
$ws.Cells.Item(90, 3).Value = "Large"
$ws.Cells.Item(90, 6).Value = "Value:             <value>"
$ws.Cells.Item(91, 2).Value = "SingleUseId200"
$ws.Cells.Item(91, 6).Value = "0"
$ws.Cells.Item(92, 2).Value = "SingleUseId201"
$ws.Cells.Item(92, 6).Value = "Parameter ID:  <value>"
$ws.Cells.Item(93, 2).Value = "SingleUseId202"
$ws.Cells.Item(93, 6).Value = "0"
$ws.Cells.Item(94, 2).Value = "SingleUseId203"
$ws.Cells.Item(94, 6).Value = "Module ID:      <value>"
$ws.Cells.Item(95, 2).Value = "SingleUseId204"
$ws.Cells.Item(95, 6).Value = "0"
$ws.Cells.Item(96, 2).Value = "SingleUseId205"
$ws.Cells.Item(96, 6).Value = "Source ID:      <value>"
$ws.Cells.Item(97, 2).Value = "SingleUseId206"
$ws.Cells.Item(97, 6).Value = "0"
$ws.Cells.Item(98, 2).Value = "SingleUseId207"
$ws.Cells.Item(98, 6).Value = "UART RX Debug"
$ws.Cells.Item(99, 2).Value = "SingleUseId210"
$ws.Cells.Item(99, 6).Value = "Type ID:          <value>"
$ws.Cells.Item(100, 2).Value = "SingleUseId211"
$ws.Cells.Item(100, 6).Value = "0"
$ws.Cells.Item(101, 2).Value = "SingleUseId212"
$ws.Cells.Item(101, 6).Value = "SignID:           <value>"
$ws.Cells.Item(102, 2).Value = "SingleUseId213"
$ws.Cells.Item(102, 6).Value = "0"
$ws.Cells.Item(103, 2).Value = "SingleUseId214"
$ws.Cells.Item(103, 6).Value = "Param2           <value>"
$ws.Cells.Item(104, 2).Value = "SingleUseId215"
$ws.Cells.Item(104, 3).Value = "Large"
$ws.Cells.Item(104, 4).Value = "Left"
$ws.Cells.Item(104, 5).Value = "LTR"
$ws.Cells.Item(104, 6).Value = "0"
$ws.Cells.Item(105, 2).Value = "SingleUseId216"
$ws.Cells.Item(105, 3).Value = "Large"
$ws.Cells.Item(105, 4).Value = "Left"
$ws.Cells.Item(105, 5).Value = "LTR"
$ws.Cells.Item(105, 6).Value = "Param1           <value>"
$ws.Cells.Item(106, 2).Value = "SingleUseId217"
$ws.Cells.Item(106, 3).Value = "Large"
$ws.Cells.Item(106, 4).Value = "Left"
$ws.Cells.Item(106, 5).Value = "LTR"
$ws.Cells.Item(106, 6).Value = "0"
$ws.Cells.Item(107, 2).Value = "SingleUseId218"
$ws.Cells.Item(107, 3).Value = "Large"
$ws.Cells.Item(107, 4).Value = "Left"
$ws.Cells.Item(107, 5).Value = "LTR"
$ws.Cells.Item(107, 6).Value = "Power             <value>"
$ws.Cells.Item(108, 2).Value = "SingleUseId219"
$ws.Cells.Item(108, 3).Value = "Large"
$ws.Cells.Item(108, 4).Value = "Left"
$ws.Cells.Item(108, 5).Value = "LTR"
$ws.Cells.Item(108, 6).Value = "0"
$ws.Cells.Item(109, 2).Value = "SingleUseId220"
$ws.Cells.Item(109, 3).Value = "Large"
$ws.Cells.Item(109, 4).Value = "Left"
$ws.Cells.Item(109, 5).Value = "LTR"
$ws.Cells.Item(109, 6).Value = "Frequency       <value>"
$ws.Cells.Item(110, 2).Value = "SingleUseId221"
$ws.Cells.Item(110, 3).Value = "Large"
$ws.Cells.Item(110, 4).Value = "Left"
$ws.Cells.Item(110, 5).Value = "LTR"
$ws.Cells.Item(110, 6).Value = "0"
$ws.Cells.Item(111, 2).Value = "SingleUseId222"
$ws.Cells.Item(111, 3).Value = "Large"
$ws.Cells.Item(111, 4).Value = "Left"
$ws.Cells.Item(111, 5).Value = "LTR"
$ws.Cells.Item(111, 6).Value = "Current           <value>"
$ws.Cells.Item(112, 2).Value = "SingleUseId223"
$ws.Cells.Item(112, 3).Value = "Large"
$ws.Cells.Item(112, 4).Value = "Left"
$ws.Cells.Item(112, 5).Value = "LTR"
$ws.Cells.Item(112, 6).Value = "0"
$ws.Cells.Item(113, 2).Value = "SingleUseId224"
$ws.Cells.Item(113, 3).Value = "Large"
$ws.Cells.Item(113, 4).Value = "Left"
$ws.Cells.Item(113, 5).Value = "LTR"
$ws.Cells.Item(113, 6).Value = "Voltage           <value>"
$ws.Cells.Item(114, 2).Value = "SingleUseId225"
$ws.Cells.Item(114, 3).Value = "Large"
$ws.Cells.Item(114, 4).Value = "Left"
$ws.Cells.Item(114, 5).Value = "LTR"
$ws.Cells.Item(114, 6).Value = "0"
$ws.Cells.Item(115, 2).Value = "SingleUseId226"
$ws.Cells.Item(115, 3).Value = "Large"
$ws.Cells.Item(115, 4).Value = "Left"
$ws.Cells.Item(115, 5).Value = "LTR"
$ws.Cells.Item(115, 6).Value = "Param2           <value>"
$ws.Cells.Item(116, 2).Value = "SingleUseId227"
$ws.Cells.Item(116, 3).Value = "Large"
$ws.Cells.Item(116, 4).Value = "Left"
$ws.Cells.Item(116, 5).Value = "LTR"
$ws.Cells.Item(116, 6).Value = "0"
$ws.Cells.Item(117, 2).Value = "SingleUseId228"
$ws.Cells.Item(117, 3).Value = "Large"
$ws.Cells.Item(117, 4).Value = "Left"
$ws.Cells.Item(117, 5).Value = "LTR"
$ws.Cells.Item(117, 6).Value = "Param1           <value>"
$ws.Cells.Item(118, 2).Value = "SingleUseId229"
$ws.Cells.Item(118, 3).Value = "Large"
$ws.Cells.Item(118, 4).Value = "Left"
$ws.Cells.Item(118, 5).Value = "LTR"
$ws.Cells.Item(118, 6).Value = "0"
$ws.Cells.Item(119, 2).Value = "SingleUseId230"
$ws.Cells.Item(119, 3).Value = "Large"
$ws.Cells.Item(119, 4).Value = "Left"
$ws.Cells.Item(119, 5).Value = "LTR"
$ws.Cells.Item(119, 6).Value = "Power             <value>"
$ws.Cells.Item(120, 2).Value = "SingleUseId231"
$ws.Cells.Item(120, 3).Value = "Large"
$ws.Cells.Item(120, 4).Value = "Left"
$ws.Cells.Item(120, 5).Value = "LTR"
$ws.Cells.Item(120, 6).Value = "0"
$ws.Cells.Item(121, 2).Value = "SingleUseId232"
$ws.Cells.Item(121, 3).Value = "Large"
$ws.Cells.Item(121, 4).Value = "Left"
$ws.Cells.Item(121, 5).Value = "LTR"
$ws.Cells.Item(121, 6).Value = "Frequency       <value>"
$ws.Cells.Item(122, 2).Value = "SingleUseId233"
$ws.Cells.Item(122, 3).Value = "Large"
$ws.Cells.Item(122, 4).Value = "Left"
$ws.Cells.Item(122, 5).Value = "LTR"
$ws.Cells.Item(122, 6).Value = "0"
$ws.Cells.Item(123, 2).Value = "SingleUseId234"
$ws.Cells.Item(123, 3).Value = "Large"
$ws.Cells.Item(123, 4).Value = "Left"
$ws.Cells.Item(123, 5).Value = "LTR"
$ws.Cells.Item(123, 6).Value = "Current           <value>"
$ws.Cells.Item(124, 2).Value = "SingleUseId235"
$ws.Cells.Item(124, 3).Value = "Large"
$ws.Cells.Item(124, 4).Value = "Left"
$ws.Cells.Item(124, 5).Value = "LTR"
$ws.Cells.Item(124, 6).Value = "0"
$ws.Cells.Item(125, 2).Value = "SingleUseId236"
$ws.Cells.Item(125, 3).Value = "Large"
$ws.Cells.Item(125, 4).Value = "Left"
$ws.Cells.Item(125, 5).Value = "LTR"
$ws.Cells.Item(125, 6).Value = "Voltage           <value>"
$ws.Cells.Item(126, 2).Value = "SingleUseId237"
$ws.Cells.Item(126, 3).Value = "Large"
$ws.Cells.Item(126, 4).Value = "Left"
$ws.Cells.Item(126, 5).Value = "LTR"
$ws.Cells.Item(126, 6).Value = "0"
$ws.Cells.Item(127, 2).Value = "SingleUseId238"
$ws.Cells.Item(127, 3).Value = "Large"
$ws.Cells.Item(127, 4).Value = "Left"
$ws.Cells.Item(127, 5).Value = "LTR"
$ws.Cells.Item(127, 6).Value = "2"
$ws.Cells.Item(128, 2).Value = "SingleUseId239"
$ws.Cells.Item(128, 3).Value = "Large"
$ws.Cells.Item(128, 4).Value = "Left"
$ws.Cells.Item(128, 5).Value = "LTR"
$ws.Cells.Item(128, 6).Value = "1"
$ws.Cells.Item(129, 2).Value = "SingleUseId240"
$ws.Cells.Item(129, 3).Value = "Large"
$ws.Cells.Item(129, 4).Value = "Left"
$ws.Cells.Item(129, 5).Value = "LTR"
$ws.Cells.Item(129, 6).Value = "3"
$ws.Cells.Item(130, 2).Value = "SingleUseId241"
$ws.Cells.Item(130, 3).Value = "Large"
$ws.Cells.Item(130, 4).Value = "Left"
$ws.Cells.Item(130, 5).Value = "LTR"
$ws.Cells.Item(130, 6).Value = "4"
$ws.Cells.Item(131, 2).Value = "SingleUseId242"
$ws.Cells.Item(131, 3).Value = "Large"
$ws.Cells.Item(131, 4).Value = "Left"
$ws.Cells.Item(131, 5).Value = "LTR"
$ws.Cells.Item(131, 6).Value = "5"
$ws.Cells.Item(132, 2).Value = "SingleUseId243"
$ws.Cells.Item(132, 3).Value = "Large"
$ws.Cells.Item(132, 4).Value = "Left"
$ws.Cells.Item(132, 5).Value = "LTR"
$ws.Cells.Item(132, 6).Value = "6"
$ws.Cells.Item(133, 2).Value = "SingleUseId244"
$ws.Cells.Item(133, 3).Value = "Large"
$ws.Cells.Item(133, 4).Value = "Left"
$ws.Cells.Item(133, 5).Value = "LTR"
$ws.Cells.Item(133, 6).Value = "7"
$ws.Cells.Item(134, 2).Value = "SingleUseId245"
$ws.Cells.Item(134, 3).Value = "Large"
$ws.Cells.Item(134, 4).Value = "Left"
$ws.Cells.Item(134, 5).Value = "LTR"
$ws.Cells.Item(134, 6).Value = "8"
$ws.Cells.Item(135, 2).Value = "SingleUseId246"
$ws.Cells.Item(135, 3).Value = "Large"
$ws.Cells.Item(135, 4).Value = "Left"
$ws.Cells.Item(135, 5).Value = "LTR"
$ws.Cells.Item(135, 6).Value = "9"
$ws.Cells.Item(136, 2).Value = "SingleUseId247"
$ws.Cells.Item(136, 3).Value = "Large"
$ws.Cells.Item(136, 4).Value = "Left"
$ws.Cells.Item(136, 5).Value = "LTR"
$ws.Cells.Item(136, 6).Value = "10"
$ws.Cells.Item(137, 2).Value = "SingleUseId248"
$ws.Cells.Item(137, 3).Value = "Large"
$ws.Cells.Item(137, 4).Value = "Left"
$ws.Cells.Item(137, 5).Value = "LTR"
$ws.Cells.Item(137, 6).Value = "9"
$ws.Cells.Item(138, 2).Value = "SingleUseId249"
$ws.Cells.Item(138, 3).Value = "Large"
$ws.Cells.Item(138, 4).Value = "Left"
$ws.Cells.Item(138, 5).Value = "LTR"
$ws.Cells.Item(138, 6).Value = "10"
$ws.Cells.Item(139, 2).Value = "SingleUseId250"
$ws.Cells.Item(139, 3).Value = "Default"
$ws.Cells.Item(139, 4).Value = "Center"
$ws.Cells.Item(139, 5).Value = "LTR"
$ws.Cells.Item(139, 6).Value = "Disable Parameter"
